$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Texas -- Bexar County): Date Published moves from 2020-07-21 to 2020-07-22
$ws.Range("B3").Value = 44034

# Row 4 (New York -- New York): scrape failed on this run -- most columns
# came back blank and the status flips to a rate-limit error message.
$blankCells = @("B4", "C4", "D4", "E4", "F4", "G4", "H4", "K4", "L4")
foreach ($cellAddr in $blankCells) {
    $ws.Range($cellAddr).Value = ""
    $ws.Range($cellAddr).Style = "Normal"
}
$ws.Range("J4").Value = $false
$ws.Range("O4").Value = "An error occurred. ... RateLimitExceededException(403, {'message': ""API rate limit exceeded for 132.145.200.60. (But here's the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)"", 'documentation_url': 'https://developer.github.com/v3/#rate-limiting'})"

# Row 39 (Delaware): error message changes from an AttributeError to an HTTPError
$ws.Range("O39").Value = "An error occurred. ... HTTPError('504 Server Error: Gateway Time-out for url: https://myhealthycommunity.dhss.delaware.gov/locations/state/')"
